# Generate Report for Handback
# The "dac72cc2-e4b4-4292-bf5d-863e2319ed89.md" file has finished localization
# handback, so its Status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" on every sheet that tracks it, and the
# "Latest Handback DateTime" timestamps for that handback are refreshed.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusHandedBack
$wsZhCn.Range("G2").Value = "2016-03-09 06:26:51"
$wsZhCn.Range("G3").Value = "2016-03-09 06:26:51"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusHandedBack
$wsDeDe.Range("G2").Value = "2016-03-09 06:27:10"
$wsDeDe.Range("G3").Value = "2016-03-09 06:27:10"
